$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Developmnet PB")

# New requirement rows describing Tree Editor equipment node support
$items = @(
    "Change Attribute to Patient attribute in the Tree editor",
    "Handle when equipment node is first node to be inserted",
    "Fix bug where when you delete a solution node, you can then insert a patient/equipment attribute node afterwards",
    "Change length of numeric attribute to be smaller"
)
$statuses = @("Not Started", "Not Started", "Done", "Not Started")
# BGR integers matching existing fills: fillId3 (orange FFC000) = 49407, fillId4 (yellow FFFF00) = 65535
$priorityColors = @(49407, 49407, 49407, 65535)

$startRow = 39
for ($i = 0; $i -lt $items.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $items[$i]
    $ws.Cells.Item($row, 1).WrapText = $true
    $ws.Cells.Item($row, 2).Interior.Color = $priorityColors[$i]
    $ws.Cells.Item($row, 3).Value = $statuses[$i]
}

$ws.Range("E38").Select()
